$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 24 and 25 hold two different species observations whose row order
# needs to be swapped. Only the columns that actually differ between the
# two rows are touched, so columns with identical values (and therefore
# untouched formatting/type) are left completely alone.
$cols = @("A","B","E","F","G","H","M","Q","R","Z","AB")

$row24 = @{}
$row25 = @{}
foreach ($col in $cols) {
    $row24[$col] = $ws.Range($col + "24").Value2
    $row25[$col] = $ws.Range($col + "25").Value2
}

foreach ($col in $cols) {
    # Row 25's original value moves up into row 24 ...
    if ($row25[$col] -eq "" -or $row25[$col] -eq $null) {
        $ws.Range($col + "24").ClearContents()
    } else {
        $ws.Range($col + "24").Value = $row25[$col]
    }

    # ... and row 24's original value moves down into row 25.
    if ($row24[$col] -eq "" -or $row24[$col] -eq $null) {
        $ws.Range($col + "25").ClearContents()
    } else {
        $ws.Range($col + "25").Value = $row24[$col]
    }
}
